$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# need NumberFormat "@" (Text) first so the literal string is preserved,
# matching the source data which keeps these as plain text cells.

# Row 2
$ws.Range("D2").Value = "43.400.13"
$ws.Range("E2").Value = "  +0.41%  "

# Row 3
$ws.Range("D3").Value = "2.409.25"
$ws.Range("E3").Value = "  +2.36%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.76"
$ws.Range("E5").Value = "  +1.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.22"
$ws.Range("E6").Value = "  +1.57%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.507"
$ws.Range("E7").Value = "  +0.39%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  -2.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.96"
$ws.Range("E10").Value = "  +2.29%  "

# Row 11
$ws.Range("E11").Value = "  +2.96%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0794"
$ws.Range("E12").Value = "  +0.69%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.48"
$ws.Range("E13").Value = "  -1.17%  "

# Row 14
$ws.Range("E14").Value = "  +2.10%  "

# Row 15
$ws.Range("D15").Value = "2.783.99"
$ws.Range("E15").Value = "  +2.30%  "

# Row 16
$ws.Range("D16").Value = "2.422.78"
$ws.Range("E16").Value = "  +3.58%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.824"
$ws.Range("E17").Value = "  +3.40%  "

# Row 18
$ws.Range("D18").Value = "43.370.86"
$ws.Range("E18").Value = "  +0.40%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.41"
$ws.Range("E19").Value = "  +2.46%  "

# Row 20
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.14"
$ws.Range("E20").Value = "  -0.56%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0900"
$ws.Range("E21").Value = "  +1.02%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.48"
$ws.Range("E22").Value = "  +0.44%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.70"
$ws.Range("E23").Value = "  +0.88%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("E24").Value = "  +0.32%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.44"
$ws.Range("E25").Value = "  +0.47%  "

# Row 26
$ws.Range("E26").Value = "  +0.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.89"
$ws.Range("E27").Value = "  +1.35%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("E28").Value = "  -4.77%  "

# Row 29
$ws.Range("E29").Value = "  +2.98%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.23"
$ws.Range("E30").Value = "  +2.52%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.11"
$ws.Range("E31").Value = "  +1.67%  "

# Row 32
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.04%  "

# Row 33
$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.31"
$ws.Range("E33").Value = "  +6.53%  "

# Row 34
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.113"
$ws.Range("E34").Value = "  +12.92%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0743"
$ws.Range("E35").Value = "  +2.51%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "132.85"
$ws.Range("E36").Value = "  +20.68%  "

# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.88"
$ws.Range("E37").Value = "  +2.32%  "

# Row 38
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.93"
$ws.Range("E38").Value = "  +6.49%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.37"
$ws.Range("E39").Value = "  -0.62%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.27"
$ws.Range("E40").Value = "  -1.73%  "

# Row 41
$ws.Range("E41").Value = "  -0.22%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.45"
$ws.Range("E42").Value = "  -4.78%  "

# Row 43
$ws.Range("D43").Value = "1.949.18"
$ws.Range("E43").Value = "  +0.25%  "

# Row 44
$ws.Range("E44").Value = "  +1.12%  "

# Row 45
$ws.Range("E45").Value = "  +1.91%  "

# Row 46
$ws.Range("E46").Value = "  +2.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.26"
$ws.Range("E47").Value = "  -1.52%  "

# Row 48
$ws.Range("D48").Value = "2.639.05"
$ws.Range("E48").Value = "  +2.15%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.48"
$ws.Range("E50").Value = "  -0.92%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.17"
$ws.Range("E51").Value = "  +0.06%  "
